$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 152.0
$ws.Range("C2").Value = 242.0
$ws.Range("D2").Value = 128.0
